$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.72%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-6.79%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.103"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.77%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07757"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-5.80%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.392"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.50%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.900"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-8.99%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.213"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.49%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.76%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9207"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.18%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1266"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.54%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1891"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.42%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08715"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.34%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03442"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.26%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09731"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.79%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.15%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006150"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.76%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.563"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.59%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3417"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-2.38%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1287"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.83%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.47%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.53%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.02120"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5,223.26%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04347"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.09%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.32%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004493"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.29%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "5.08%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02161"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.41%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04920"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-5.68%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007696"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.35%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.32%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1336"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.18%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002007"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.22%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008863"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.24%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006847"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.66%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000755"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "1.16%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003018"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "3.16%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001310"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-22.17%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002115"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "1.16%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002014"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "1.16%"
